$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates replaced with an invalid negative serial number
$ws.Range("B6").Value = -692739
$ws.Range("B12").Value = -692739
$ws.Range("F12").Value = -692739

# Traveler name / identification
$ws.Range("B7").Value = "John Doe"
$ws.Range("F8").Value = "CC123"

# Destination city
$ws.Range("B11").Value = "Nueva York"

# Trip reason
$ws.Range("B13").Value = "Asistir a una conferencia de ventas"

# Expense row 20
$ws.Range("B20").Value = "'12345"
$ws.Range("C20").Value = "Aerolínea XYZ"
$ws.Range("D20").Value = "'1234567890"
$ws.Range("E20").Value = "Tiquete aéreo ida y vuelta"
$ws.Range("G20").Value = 1000

# Advance amount row 32
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 1000

# Signatures
$ws.Range("A40").Value = "Bob Johnson"
$ws.Range("E40").Value = "Jane Smith"

# Bank
$ws.Range("B45").Value = "Banco XYZ"
